# Update countries & provincias Spain
# - Refresh the "last updated" timestamp in A1
# - Update case counts for several provinces (Ciudad Real, Albacete, Zaragoza,
#   Toledo, Cuenca, Huesca, Teruel)
# - Guadalajara's totals overtook Cordoba/Jaen, so it moves up above them in
#   the (descending, by "Casos totales") sorted list; Cordoba and Jaen shift
#   down one row each, keeping their own data intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 12:22"

# --- Straight numeric refreshes (no reordering) ------------------------
# Ciudad Real (row 8)
$ws.Range("B8").Value = 5563
$ws.Range("C8").Value = 2943
$ws.Range("D8").Value = 9672
$ws.Range("E8").Value = 647

# Albacete (row 11)
$ws.Range("B11").Value = 3543
$ws.Range("C11").Value = 2943
$ws.Range("D11").Value = 9672
$ws.Range("E11").Value = 331

# Zaragoza (row 14)
$ws.Range("B14").Value = 3186
$ws.Range("C14").Value = 722
$ws.Range("D14").Value = 2094
$ws.Range("E14").Value = 370

# Toledo (row 15)
$ws.Range("B15").Value = 3098
$ws.Range("C15").Value = 2943
$ws.Range("D15").Value = 9672
$ws.Range("E15").Value = 442

# --- Guadalajara moves above Cordoba/Jaen (rows 34-36) ------------------
# Row 34: was Cordoba -> now Guadalajara, with refreshed totals
$ws.Range("A34").Value = "Guadalajara"
$ws.Range("B34").Value = 1195
$ws.Range("C34").Value = 2943
$ws.Range("D34").Value = 9672
$ws.Range("E34").Value = 154

# Row 35: was Jaen -> now Cordoba (its own unchanged figures)
$ws.Range("A35").Value = "Cordoba"
$ws.Range("B35").Value = 1169
$ws.Range("C35").Value = 197
$ws.Range("D35").Value = 916
$ws.Range("E35").Value = 56

# Row 36: was Guadalajara -> now Jaen (its own unchanged figures)
$ws.Range("A36").Value = "Jaen"
$ws.Range("B36").Value = 1169
$ws.Range("C36").Value = 205
$ws.Range("D36").Value = 857
$ws.Range("E36").Value = 107

# Cuenca (row 41)
$ws.Range("B41").Value = 930
$ws.Range("C41").Value = 2943
$ws.Range("D41").Value = 9672
$ws.Range("E41").Value = 140

# Huesca (row 47)
$ws.Range("B47").Value = 508
$ws.Range("C47").Value = 113
$ws.Range("D47").Value = 326
$ws.Range("E47").Value = 69

# Teruel (row 48)
$ws.Range("B48").Value = 502
$ws.Range("C48").Value = 118
$ws.Range("D48").Value = 332
$ws.Range("E48").Value = 52
